# Fixed DropDown issue Commit
# - Update the stale verification email (shared by SignIn!A2 and CreateAccount!D2)
# - Add a hyperlink on CreateAccount!D2 for the updated email
# - Make CreateAccount the active tab / selected sheet (instead of SignIn)
# - Update the remembered cell selections on both sheets

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # SignIn
$ws2 = $wb.Worksheets.Item(2)  # CreateAccount

# --- Update the email text on both sheets so the shared string is edited in place ---
$ws1.Range("A2").Value = "testjaga18042018@gmail.com"
$ws2.Range("D2").Value = "testjaga18042018@gmail.com"

# --- Add hyperlink on CreateAccount!D2 pointing at the updated email ---
$null = $ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:testjaga18042018@gmail.com")
# Restore the existing "Hyperlink" cell style (Add() stamps a fresh style record,
# this keeps the cell's style identical to before the hyperlink was added)
$ws2.Range("D2").Style = "Hyperlink"

# --- Update the saved selections on each sheet ---
$null = $ws1.Range("A2").Select()
$null = $ws2.Range("F3").Select()

# --- Make CreateAccount the active/selected tab ---
$null = $ws2.Activate()
